# The commit swaps the deck's theme palette from the custom "Integral"
# colour set over to the stock "Office Theme" colour set (the twelve
# named theme colours: dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# We reach the shared theme part through any slide's ThemeColorScheme -
# all slides share the single master theme, so this updates the one
# <a:clrScheme> used across the whole presentation.

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Colors(1).RGB  = RGB(0, 0, 0)         # dk1
$tcs.Colors(2).RGB  = RGB(255, 255, 255)   # lt1
$tcs.Colors(3).RGB  = RGB(68, 84, 106)     # dk2
$tcs.Colors(4).RGB  = RGB(231, 230, 230)   # lt2
$tcs.Colors(5).RGB  = RGB(91, 155, 213)    # accent1
$tcs.Colors(6).RGB  = RGB(237, 125, 49)    # accent2
$tcs.Colors(7).RGB  = RGB(165, 165, 165)   # accent3
$tcs.Colors(8).RGB  = RGB(255, 192, 0)     # accent4
$tcs.Colors(9).RGB  = RGB(68, 114, 196)    # accent5
$tcs.Colors(10).RGB = RGB(112, 173, 71)    # accent6
$tcs.Colors(11).RGB = RGB(5, 99, 193)      # hlink
$tcs.Colors(12).RGB = RGB(149, 79, 114)    # folHlink
